$p = $ppt.ActivePresentation
Write-Host "TemplateName:" $p.TemplateName
try { $p.TemplateName = "Office Theme" ; Write-Host "set ok" } catch { Write-Host "ERR:" $_.Exception.Message }
Write-Host "TemplateName2:" $p.TemplateName
